$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update marking score (points per correct answer) and the resulting total,
# plus the "correct/total marks" summary cell.
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 30
$ws.Range("E12").Value = "30/140"
